# Apply the "Add files via upload" revision to WS_holdings:
#  1. Bump the confidentiality footer's "as of" date 2021-06-09 -> 2021-06-10.
#  2. Refresh the Weight (D) / Percent Change (E) figures for rows 2-13.
# The worksheet ships protected (sheet-level protection with no editable
# cells), so we unprotect, write the values, then restore protection
# afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect()
}

# --- Footer text: update the "as of" date -------------------------------
$footer = $ws.Range("A16").Value2
$ws.Range("A16").Value2 = $footer -replace "2021-06-09", "2021-06-10"

# --- Weight / Percent Change refresh, rows 2-13 --------------------------
$ws.Range("D2").Value2 = 0.02817013739599431
$ws.Range("E2").Value2 = 0.006702840727736969

$ws.Range("D3").Value2 = 0.02231158116676947
$ws.Range("E3").Value2 = 0.006829957607159765

$ws.Range("D4").Value2 = 0.05874787504580373
$ws.Range("E4").Value2 = 0.005116279069767415

$ws.Range("D5").Value2 = 0.1379121255100299
$ws.Range("E5").Value2 = 0.01753539253539249

$ws.Range("D6").Value2 = 0.02199466259615863
$ws.Range("E6").Value2 = -0.00197498354180381

$ws.Range("D7").Value2 = 0.1264182806665426
$ws.Range("E7").Value2 = 0.001274697259400881

$ws.Range("D8").Value2 = 0.09218056558170906
$ws.Range("E8").Value2 = -0.005472455308281665

$ws.Range("D9").Value2 = 0.03174136828950616
$ws.Range("E9").Value2 = -0.005974454058508383

$ws.Range("D10").Value2 = 0.1094342945128714
$ws.Range("E10").Value2 = -0.01132565911622718

$ws.Range("D11").Value2 = 0.2819016375371714
$ws.Range("E11").Value2 = 0.008719101123595419

$ws.Range("D12").Value2 = 0.08918747169744343
$ws.Range("E12").Value2 = 0.003205128205128194

$ws.Range("E13").Value2 = 0.003988106095089039

# --- Restore sheet protection --------------------------------------------
if ($wasProtected) {
    $ws.Protect()
}
